$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 2.56
$ws.Range("T2").Value = 1.66
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 28
$ws.Range("K3").Value = 17.5
$ws.Range("L3").Value = 1.17
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 7.6
$ws.Range("O3").Value = 1.11
$ws.Range("R3").Value = 2.02
$ws.Range("S3").Value = 1.71
$ws.Range("T3").Value = 1.95
$ws.Range("U3").Value = 1.68
$ws.Range("V3").Value = 5.5
$ws.Range("W3").Value = 1.03
$ws.Range("X3").Value = 60
$ws.Range("Y3").Value = 17
$ws.Range("Z3").Value = 12
$ws.Range("AA3").Value = 11
$ws.Range("AB3").Value = 90
$ws.Range("AC3").Value = 27
$ws.Range("AD3").Value = 16.5
$ws.Range("AE3").Value = 16.5
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 85
$ws.Range("AH3").Value = 48
$ws.Range("AI3").Value = 48
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 2.98
$ws.Range("G4").Value = 3.15
$ws.Range("H4").Value = 2.54
$ws.Range("U6").Value = 2.2
$ws.Range("X6").Value = 15
$ws.Range("AC6").Value = 8.199999999999999
$ws.Range("AF6").Value = 29
$ws.Range("AI6").Value = 42
$ws.Range("AO6").Value = 19
$ws.Range("G7").Value = 2.72
$ws.Range("I7").Value = 3.6
$ws.Range("P8").Value = 1.97
$ws.Range("Q8").Value = 1.71
$ws.Range("Q9").Value = 1.68
$ws.Range("F11").Value = 1.63
$ws.Range("G11").Value = 1.69
$ws.Range("H11").Value = 5
$ws.Range("J11").Value = 4.2
$ws.Range("K11").Value = 5.4
$ws.Range("P11").Value = 2.28
$ws.Range("I12").Value = 3.85
$ws.Range("Q12").Value = 1.72
$ws.Range("X12").Value = 21
$ws.Range("AA12").Value = 75
$ws.Range("AB12").Value = 12
$ws.Range("AC12").Value = 9
$ws.Range("AE12").Value = 42
$ws.Range("AF12").Value = 15
$ws.Range("AH12").Value = 16
$ws.Range("AI12").Value = 46
$ws.Range("AK12").Value = 19.5
$ws.Range("AL12").Value = 32
$ws.Range("AM12").Value = 75
$ws.Range("AO12").Value = 34
$ws.Range("F13").Value = 6.2
$ws.Range("H13").Value = 1.56
$ws.Range("I13").Value = 1.58
$ws.Range("Z13").Value = 12
$ws.Range("AG13").Value = 26
$ws.Range("AI13").Value = 28
$ws.Range("AL13").Value = 70
$ws.Range("AM13").Value = 80
$ws.Range("AO13").Value = 5.9
$ws.Range("F16").Value = 3.65
$ws.Range("G16").Value = 4.8
$ws.Range("H16").Value = 2.14
$ws.Range("I16").Value = 2.42
$ws.Range("K16").Value = 3.5
$ws.Range("P16").Value = 1.57
$ws.Range("Q16").Value = 2.42

$wb.Save()
